$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its original text formatting so that
# numeric-looking values (e.g. "26.16") are stored as text, matching the
# source data which used inline strings rather than numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.021.84"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").Value = "1.910.19"
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "0.7859"
$ws.Range("E5").Value = "  +5.27%  "
$ws.Range("D6").Value = "241.69"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").Value = "0.3155"
$ws.Range("E8").Value = "  +2.52%  "
$ws.Range("D9").Value = "26.16"
$ws.Range("E9").Value = "  +0.51%  "
$ws.Range("D10").Value = "0.06902"
$ws.Range("E10").Value = "  -0.08%  "
$ws.Range("D11").Value = "0.07965"
$ws.Range("E11").Value = "  -0.60%  "
$ws.Range("D12").Value = "1.908.06"
$ws.Range("E12").Value = "  +0.06%  "
$ws.Range("D13").Value = "0.7430"
$ws.Range("E13").Value = "  -2.30%  "
$ws.Range("D14").Value = "5.211"
$ws.Range("E14").Value = "  -0.86%  "
$ws.Range("D15").Value = "92.96"
$ws.Range("E15").Value = "  +1.73%  "
$ws.Range("D16").Value = "30.029.22"
$ws.Range("E16").Value = "  -0.12%  "
$ws.Range("D17").Value = "13.98"
$ws.Range("E17").Value = "  -0.77%  "
$ws.Range("D18").Value = "5.883"
$ws.Range("E18").Value = "  -4.94%  "
$ws.Range("D19").Value = "246.35"
$ws.Range("E19").Value = "  +3.87%  "
$ws.Range("D20").Value = "0.000007755"
$ws.Range("E20").Value = "  -0.21%  "
$ws.Range("E21").Value = "  +0.19%  "
$ws.Range("D22").Value = "2.150.27"
$ws.Range("E22").Value = "  -0.29%  "
$ws.Range("E23").Value = "  +0.29%  "
$ws.Range("D24").Value = "6.878"
$ws.Range("E24").Value = "  -3.42%  "
$ws.Range("D25").Value = "169.38"
$ws.Range("E25").Value = "  +1.31%  "
$ws.Range("D26").Value = "9.288"
$ws.Range("E26").Value = "  -0.66%  "
$ws.Range("D27").Value = "0.1375"
$ws.Range("E27").Value = "  +8.41%  "
$ws.Range("D28").Value = "18.94"
$ws.Range("E28").Value = "  +0.42%  "
$ws.Range("D29").Value = "2.030"
$ws.Range("E29").Value = "  -0.95%  "
$ws.Range("D30").Value = "1.378"
$ws.Range("E30").Value = "  +1.47%  "
$ws.Range("D31").Value = "1.520"
$ws.Range("E31").Value = "  -0.61%  "
$ws.Range("D32").Value = "4.325"
$ws.Range("E32").Value = "  +0.53%  "
$ws.Range("D33").Value = "4.083"
$ws.Range("E33").Value = "  +0.93%  "
$ws.Range("D34").Value = "0.05501"
$ws.Range("E34").Value = "  +2.67%  "
$ws.Range("D35").Value = "1.257"
$ws.Range("E35").Value = "  -2.61%  "
$ws.Range("D36").Value = "0.7340"
$ws.Range("E36").Value = "  -1.18%  "
$ws.Range("E37").Value = "  +0.52%  "
$ws.Range("D38").Value = "0.01935"
$ws.Range("E38").Value = "  -0.45%  "
$ws.Range("D39").Value = "2.792"
$ws.Range("E39").Value = "  +1.18%  "
$ws.Range("D40").Value = "6.130"
$ws.Range("E40").Value = "  -2.03%  "
$ws.Range("D41").Value = "0.4420"
$ws.Range("E41").Value = "  -1.06%  "
$ws.Range("D42").Value = "72.10"
$ws.Range("E42").Value = "  -1.37%  "
$ws.Range("E43").Value = "  +0.25%  "
$ws.Range("D44").Value = "0.8372"
$ws.Range("E44").Value = "  +0.79%  "
$ws.Range("D45").Value = "1.878"
$ws.Range("E45").Value = "  -4.57%  "
$ws.Range("D46").Value = "100.44"
$ws.Range("E46").Value = "  -0.99%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "9.800"
$ws.Range("E47").Value = "  -0.22%  "
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "7.529"
$ws.Range("E48").Value = "  -2.40%  "
$ws.Range("D49").Value = "980.84"
$ws.Range("E49").Value = "  +8.48%  "
$ws.Range("D50").Value = "2.057.90"
$ws.Range("E50").Value = "  -0.35%  "
$ws.Range("E51").Value = "  -1.29%  "
